$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 with row 4, and row 3 with row 5
# (dates / variety / quality / volume / price-unit-label values trade places).

# --- Save original row 2 and row 4 values ---
$row2 = $ws.Range("A2:T2").Value2
$row4 = $ws.Range("A4:T4").Value2

$ws.Range("A2:T2").Value2 = $row4
$ws.Range("A4:T4").Value2 = $row2

# --- Save original row 3 and row 5 values ---
$row3 = $ws.Range("A3:T3").Value2
$row5 = $ws.Range("A5:T5").Value2

$ws.Range("A3:T3").Value2 = $row5
$ws.Range("A5:T5").Value2 = $row3
